$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("__data")

$ws.Range("G6").Value = "damage:+24|pullStrength:+20"
$ws.Range("G7").Value = "projectileSpeed:+18|split:+1|crit:+6"
$ws.Range("F8").Value = "relic:30050002|skill:70030002"
$ws.Range("G8").Value = "slow:+12%|shield:+30|duration:+2"
$ws.Range("F9").Value = "weapon:20020003|skill:70040002"
$ws.Range("G9").Value = "damageMultiplier:+12|ricochet:+1"

$wb.Save()
